$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.604.01'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '3.362.62'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '557.51'
$ws.Range('E5').Value = '  -3.36%  '
$ws.Range('D6').Value = '175.96'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  -1.29%  '
$ws.Range('D8').Value = '3.355.33'
$ws.Range('E8').Value = '  -2.74%  '
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').Value = '0.163'
$ws.Range('E11').Value = '  +2.27%  '
$ws.Range('D12').Value = '54.59'
$ws.Range('E12').Value = '  -2.02%  '
$ws.Range('D13').Value = '0.0000274'
$ws.Range('E13').Value = '  -0.84%  '
$ws.Range('D14').Value = '9.08'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').Value = '3.891.65'
$ws.Range('E15').Value = '  -2.93%  '
$ws.Range('D16').Value = '18.43'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').Value = '3.357.45'
$ws.Range('E18').Value = '  -2.91%  '
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('D20').Value = '64.424.44'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('D21').Value = '0.985'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').Value = '458.06'
$ws.Range('E22').Value = '  +12.01%  '
$ws.Range('D23').Value = '4.85'
$ws.Range('E23').Value = '  +11.00%  '
$ws.Range('D24').Value = '4.10'
$ws.Range('E24').Value = '  -3.42%  '
$ws.Range('D25').Value = '85.45'
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('D26').Value = '13.36'
$ws.Range('E26').Value = '  -1.39%  '
$ws.Range('D27').Value = '10.81'
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('D28').Value = '2.85'
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('D29').Value = '8.80'
$ws.Range('E29').Value = '  -3.33%  '
$ws.Range('D30').Value = '29.96'
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').Value = '6.60'
$ws.Range('E31').Value = '  -1.11%  '
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('D33').Value = '580.60'
$ws.Range('E33').Value = '  -1.33%  '
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('D35').Value = '58.63'
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('E37').Value = '  -8.72%  '
$ws.Range('D38').Value = '3.53'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('D39').Value = '35.78'
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('D40').Value = '0.0₃0756'
$ws.Range('E40').Value = '  -3.60%  '
$ws.Range('D41').Value = '0.369'
$ws.Range('E41').Value = '  -2.94%  '
$ws.Range('D42').Value = '3.101.07'
$ws.Range('E42').Value = '  -3.39%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '2.81'
$ws.Range('E44').Value = '  -5.25%  '
$ws.Range('D45').Value = '3.23'
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('E47').Value = '  -2.06%  '
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('D49').Value = '2.59'
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('D50').Value = '8.29'
$ws.Range('E50').Value = '  -2.96%  '
$ws.Range('D51').Value = '134.94'
$ws.Range('E51').Value = '  -2.05%  '
